$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "24.859.32", "  -4.23%  ", 0)
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.632.21", "  -6.52%  ", 0)
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.9986", "  -0.14%  ", 1)
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "232.21", "  -6.66%  ", 1)
    ,@(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.000", "  +0.03%  ", 1)
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4711", "  -6.73%  ", 1)
    ,@(8, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "38.90", "  -4.17%  ", 1)
    ,@(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2546", "  -7.54%  ", 1)
    ,@(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06053", "  -2.19%  ", 1)
    ,@(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06975", "  -4.06%  ", 1)
    ,@(12, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.646.61", "  -5.80%  ", 0)
    ,@(13, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "14.48", "  -4.67%  ", 1)
    ,@(14, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.5927", "  -9.53%  ", 1)
    ,@(15, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.302", "  -7.62%  ", 1)
    ,@(16, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "72.90", "  -6.23%  ", 1)
    ,@(17, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.000", "  -0.02%  ", 1)
    ,@(18, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "0.9993", "  -0.03%  ", 1)
    ,@(19, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "24.867.58", "  -4.27%  ", 0)
    ,@(20, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000006526", "  -4.69%  ", 1)
    ,@(21, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "11.08", "  -6.64%  ", 1)
    ,@(22, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.851.99", "  -5.93%  ", 0)
    ,@(23, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.307", "  -3.36%  ", 1)
    ,@(24, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.506", "  -2.58%  ", 1)
    ,@(25, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "5.193", "  -3.85%  ", 1)
    ,@(26, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "133.20", "  -2.77%  ", 1)
    ,@(27, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "14.73", "  -3.44%  ", 1)
    ,@(28, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.377", "  -8.64%  ", 1)
    ,@(29, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "103.13", "  -2.58%  ", 1)
    ,@(30, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.621", "  -9.06%  ", 1)
    ,@(31, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "3.782", "  -2.46%  ", 1)
    ,@(32, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.07668", "  -6.50%  ", 1)
    ,@(33, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.504", "  -4.03%  ", 1)
    ,@(34, "Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "0.9993", "  +0.02%  ", 1)
    ,@(35, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.04271", "  -8.87%  ", 1)
    ,@(36, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.580", "  -2.76%  ", 1)
    ,@(37, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.9149", "  -8.22%  ", 1)
    ,@(38, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.5743", "  -7.08%  ", 1)
    ,@(39, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.527", "  -8.10%  ", 1)
    ,@(40, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01524", "  -5.60%  ", 1)
    ,@(41, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "0.9985", "  -0.13%  ", 1)
    ,@(42, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.8086", "  +6.06%  ", 1)
    ,@(43, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "96.78", "  -4.12%  ", 1)
    ,@(44, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.749", "  -9.20%  ", 1)
    ,@(45, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.3654", "  -7.10%  ", 1)
    ,@(46, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "4.691", "  -6.40%  ", 1)
    ,@(47, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.05195", "  -1.92%  ", 1)
    ,@(48, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1084", "  -5.82%  ", 1)
    ,@(49, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "5.974", "  -5.75%  ", 1)
    ,@(50, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "29.22", "  -4.92%  ", 1)
    ,@(51, "TrueUSD", "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd", "0.9993", "  -0.30%  ", 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $dCell = $ws.Cells.Item($r, 4)
    if ($row[5] -eq 1) {
        $dCell.NumberFormat = "@"
        $dCell.Value = $row[3]
        $dCell.ClearFormats()
    } else {
        $dCell.Value = $row[3]
    }
    $ws.Cells.Item($r, 5).Value = $row[4]
}

